$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date column C for all data rows from 45184 to 45186
$ws.Range("C2:C301").Value = 45186

# Add a friendly-name second argument to each HYPERLINK formula, using the "Beteckning" (column A) text of the row
$ws.Cells.Item(2, 19).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/artfynd/A 64675-2021.xlsx", "A 64675-2021")'
$ws.Cells.Item(2, 20).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/kartor/A 64675-2021.png", "A 64675-2021")'
$ws.Cells.Item(2, 22).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/klagomål/A 64675-2021.docx", "A 64675-2021")'
$ws.Cells.Item(2, 23).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/klagomålsmail/A 64675-2021.docx", "A 64675-2021")'
$ws.Cells.Item(2, 24).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/tillsyn/A 64675-2021.docx", "A 64675-2021")'
$ws.Cells.Item(2, 25).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/tillsynsmail/A 64675-2021.docx", "A 64675-2021")'
$ws.Cells.Item(3, 19).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/artfynd/A 55504-2022.xlsx", "A 55504-2022")'
$ws.Cells.Item(3, 20).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/kartor/A 55504-2022.png", "A 55504-2022")'
$ws.Cells.Item(3, 21).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/knärot/A 55504-2022.png", "A 55504-2022")'
$ws.Cells.Item(3, 22).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/klagomål/A 55504-2022.docx", "A 55504-2022")'
$ws.Cells.Item(3, 23).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/klagomålsmail/A 55504-2022.docx", "A 55504-2022")'
$ws.Cells.Item(3, 24).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/tillsyn/A 55504-2022.docx", "A 55504-2022")'
$ws.Cells.Item(3, 25).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/tillsynsmail/A 55504-2022.docx", "A 55504-2022")'
$ws.Cells.Item(4, 19).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/artfynd/A 73613-2021.xlsx", "A 73613-2021")'
$ws.Cells.Item(4, 20).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/kartor/A 73613-2021.png", "A 73613-2021")'
$ws.Cells.Item(4, 22).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/klagomål/A 73613-2021.docx", "A 73613-2021")'
$ws.Cells.Item(4, 23).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/klagomålsmail/A 73613-2021.docx", "A 73613-2021")'
$ws.Cells.Item(4, 24).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/tillsyn/A 73613-2021.docx", "A 73613-2021")'
$ws.Cells.Item(4, 25).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/tillsynsmail/A 73613-2021.docx", "A 73613-2021")'
$ws.Cells.Item(5, 19).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/artfynd/A 59728-2021.xlsx", "A 59728-2021")'
$ws.Cells.Item(5, 20).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/kartor/A 59728-2021.png", "A 59728-2021")'
$ws.Cells.Item(5, 21).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/knärot/A 59728-2021.png", "A 59728-2021")'
$ws.Cells.Item(5, 22).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/klagomål/A 59728-2021.docx", "A 59728-2021")'
$ws.Cells.Item(5, 23).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/klagomålsmail/A 59728-2021.docx", "A 59728-2021")'
$ws.Cells.Item(5, 24).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/tillsyn/A 59728-2021.docx", "A 59728-2021")'
$ws.Cells.Item(5, 25).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/tillsynsmail/A 59728-2021.docx", "A 59728-2021")'
$ws.Cells.Item(6, 19).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/artfynd/A 22181-2023.xlsx", "A 22181-2023")'
$ws.Cells.Item(6, 20).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/kartor/A 22181-2023.png", "A 22181-2023")'
$ws.Cells.Item(6, 21).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/knärot/A 22181-2023.png", "A 22181-2023")'
$ws.Cells.Item(6, 22).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/klagomål/A 22181-2023.docx", "A 22181-2023")'
$ws.Cells.Item(6, 23).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/klagomålsmail/A 22181-2023.docx", "A 22181-2023")'
$ws.Cells.Item(6, 24).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/tillsyn/A 22181-2023.docx", "A 22181-2023")'
$ws.Cells.Item(6, 25).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/tillsynsmail/A 22181-2023.docx", "A 22181-2023")'
$ws.Cells.Item(7, 19).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/artfynd/A 13707-2021.xlsx", "A 13707-2021")'
$ws.Cells.Item(7, 20).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/kartor/A 13707-2021.png", "A 13707-2021")'
$ws.Cells.Item(7, 22).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/klagomål/A 13707-2021.docx", "A 13707-2021")'
$ws.Cells.Item(7, 23).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/klagomålsmail/A 13707-2021.docx", "A 13707-2021")'
$ws.Cells.Item(7, 24).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/tillsyn/A 13707-2021.docx", "A 13707-2021")'
$ws.Cells.Item(7, 25).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/tillsynsmail/A 13707-2021.docx", "A 13707-2021")'
$ws.Cells.Item(8, 19).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/artfynd/A 40633-2022.xlsx", "A 40633-2022")'
$ws.Cells.Item(8, 20).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/kartor/A 40633-2022.png", "A 40633-2022")'
$ws.Cells.Item(8, 21).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/knärot/A 40633-2022.png", "A 40633-2022")'
$ws.Cells.Item(8, 22).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/klagomål/A 40633-2022.docx", "A 40633-2022")'
$ws.Cells.Item(8, 23).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/klagomålsmail/A 40633-2022.docx", "A 40633-2022")'
$ws.Cells.Item(8, 24).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/tillsyn/A 40633-2022.docx", "A 40633-2022")'
$ws.Cells.Item(8, 25).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/tillsynsmail/A 40633-2022.docx", "A 40633-2022")'
$ws.Cells.Item(9, 19).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/artfynd/A 16246-2023.xlsx", "A 16246-2023")'
$ws.Cells.Item(9, 20).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/kartor/A 16246-2023.png", "A 16246-2023")'
$ws.Cells.Item(9, 21).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/knärot/A 16246-2023.png", "A 16246-2023")'
$ws.Cells.Item(9, 22).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/klagomål/A 16246-2023.docx", "A 16246-2023")'
$ws.Cells.Item(9, 23).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/klagomålsmail/A 16246-2023.docx", "A 16246-2023")'
$ws.Cells.Item(9, 24).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/tillsyn/A 16246-2023.docx", "A 16246-2023")'
$ws.Cells.Item(9, 25).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/tillsynsmail/A 16246-2023.docx", "A 16246-2023")'
$ws.Cells.Item(10, 19).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/artfynd/A 16084-2019.xlsx", "A 16084-2019")'
$ws.Cells.Item(10, 20).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/kartor/A 16084-2019.png", "A 16084-2019")'
$ws.Cells.Item(10, 22).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/klagomål/A 16084-2019.docx", "A 16084-2019")'
$ws.Cells.Item(10, 23).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/klagomålsmail/A 16084-2019.docx", "A 16084-2019")'
$ws.Cells.Item(10, 24).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/tillsyn/A 16084-2019.docx", "A 16084-2019")'
$ws.Cells.Item(10, 25).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/tillsynsmail/A 16084-2019.docx", "A 16084-2019")'
$ws.Cells.Item(11, 19).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/artfynd/A 66569-2019.xlsx", "A 66569-2019")'
$ws.Cells.Item(11, 20).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/kartor/A 66569-2019.png", "A 66569-2019")'
$ws.Cells.Item(11, 22).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/klagomål/A 66569-2019.docx", "A 66569-2019")'
$ws.Cells.Item(11, 23).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/klagomålsmail/A 66569-2019.docx", "A 66569-2019")'
$ws.Cells.Item(11, 24).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/tillsyn/A 66569-2019.docx", "A 66569-2019")'
$ws.Cells.Item(11, 25).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/tillsynsmail/A 66569-2019.docx", "A 66569-2019")'
$ws.Cells.Item(12, 19).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/artfynd/A 6094-2022.xlsx", "A 6094-2022")'
$ws.Cells.Item(12, 20).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/kartor/A 6094-2022.png", "A 6094-2022")'
$ws.Cells.Item(12, 22).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/klagomål/A 6094-2022.docx", "A 6094-2022")'
$ws.Cells.Item(12, 23).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/klagomålsmail/A 6094-2022.docx", "A 6094-2022")'
$ws.Cells.Item(12, 24).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/tillsyn/A 6094-2022.docx", "A 6094-2022")'
$ws.Cells.Item(12, 25).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/tillsynsmail/A 6094-2022.docx", "A 6094-2022")'
$ws.Cells.Item(13, 19).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/artfynd/A 22913-2023.xlsx", "A 22913-2023")'
$ws.Cells.Item(13, 20).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/kartor/A 22913-2023.png", "A 22913-2023")'
$ws.Cells.Item(13, 21).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/knärot/A 22913-2023.png", "A 22913-2023")'
$ws.Cells.Item(13, 22).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/klagomål/A 22913-2023.docx", "A 22913-2023")'
$ws.Cells.Item(13, 23).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/klagomålsmail/A 22913-2023.docx", "A 22913-2023")'
$ws.Cells.Item(13, 24).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/tillsyn/A 22913-2023.docx", "A 22913-2023")'
$ws.Cells.Item(13, 25).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/tillsynsmail/A 22913-2023.docx", "A 22913-2023")'
$ws.Cells.Item(14, 19).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/artfynd/A 66895-2018.xlsx", "A 66895-2018")'
$ws.Cells.Item(14, 20).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/kartor/A 66895-2018.png", "A 66895-2018")'
$ws.Cells.Item(14, 22).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/klagomål/A 66895-2018.docx", "A 66895-2018")'
$ws.Cells.Item(14, 23).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/klagomålsmail/A 66895-2018.docx", "A 66895-2018")'
$ws.Cells.Item(14, 24).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/tillsyn/A 66895-2018.docx", "A 66895-2018")'
$ws.Cells.Item(14, 25).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/tillsynsmail/A 66895-2018.docx", "A 66895-2018")'
$ws.Cells.Item(15, 19).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/artfynd/A 67989-2019.xlsx", "A 67989-2019")'
$ws.Cells.Item(15, 20).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/kartor/A 67989-2019.png", "A 67989-2019")'
$ws.Cells.Item(15, 22).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/klagomål/A 67989-2019.docx", "A 67989-2019")'
$ws.Cells.Item(15, 23).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/klagomålsmail/A 67989-2019.docx", "A 67989-2019")'
$ws.Cells.Item(15, 24).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/tillsyn/A 67989-2019.docx", "A 67989-2019")'
$ws.Cells.Item(15, 25).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/tillsynsmail/A 67989-2019.docx", "A 67989-2019")'
$ws.Cells.Item(16, 19).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/artfynd/A 40578-2021.xlsx", "A 40578-2021")'
$ws.Cells.Item(16, 20).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/kartor/A 40578-2021.png", "A 40578-2021")'
$ws.Cells.Item(16, 21).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/knärot/A 40578-2021.png", "A 40578-2021")'
$ws.Cells.Item(16, 22).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/klagomål/A 40578-2021.docx", "A 40578-2021")'
$ws.Cells.Item(16, 23).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/klagomålsmail/A 40578-2021.docx", "A 40578-2021")'
$ws.Cells.Item(16, 24).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/tillsyn/A 40578-2021.docx", "A 40578-2021")'
$ws.Cells.Item(16, 25).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/tillsynsmail/A 40578-2021.docx", "A 40578-2021")'
$ws.Cells.Item(17, 19).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/artfynd/A 48937-2021.xlsx", "A 48937-2021")'
$ws.Cells.Item(17, 20).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/kartor/A 48937-2021.png", "A 48937-2021")'
$ws.Cells.Item(17, 21).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/knärot/A 48937-2021.png", "A 48937-2021")'
$ws.Cells.Item(17, 22).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/klagomål/A 48937-2021.docx", "A 48937-2021")'
$ws.Cells.Item(17, 23).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/klagomålsmail/A 48937-2021.docx", "A 48937-2021")'
$ws.Cells.Item(17, 24).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/tillsyn/A 48937-2021.docx", "A 48937-2021")'
$ws.Cells.Item(17, 25).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/tillsynsmail/A 48937-2021.docx", "A 48937-2021")'
$ws.Cells.Item(18, 19).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/artfynd/A 51992-2021.xlsx", "A 51992-2021")'
$ws.Cells.Item(18, 20).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/kartor/A 51992-2021.png", "A 51992-2021")'
$ws.Cells.Item(18, 22).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/klagomål/A 51992-2021.docx", "A 51992-2021")'
$ws.Cells.Item(18, 23).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/klagomålsmail/A 51992-2021.docx", "A 51992-2021")'
$ws.Cells.Item(18, 24).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/tillsyn/A 51992-2021.docx", "A 51992-2021")'
$ws.Cells.Item(18, 25).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/tillsynsmail/A 51992-2021.docx", "A 51992-2021")'
$ws.Cells.Item(19, 19).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/artfynd/A 63118-2021.xlsx", "A 63118-2021")'
$ws.Cells.Item(19, 20).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/kartor/A 63118-2021.png", "A 63118-2021")'
$ws.Cells.Item(19, 21).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/knärot/A 63118-2021.png", "A 63118-2021")'
$ws.Cells.Item(19, 22).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/klagomål/A 63118-2021.docx", "A 63118-2021")'
$ws.Cells.Item(19, 23).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/klagomålsmail/A 63118-2021.docx", "A 63118-2021")'
$ws.Cells.Item(19, 24).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/tillsyn/A 63118-2021.docx", "A 63118-2021")'
$ws.Cells.Item(19, 25).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/tillsynsmail/A 63118-2021.docx", "A 63118-2021")'
$ws.Cells.Item(20, 19).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/artfynd/A 52510-2022.xlsx", "A 52510-2022")'
$ws.Cells.Item(20, 20).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/kartor/A 52510-2022.png", "A 52510-2022")'
$ws.Cells.Item(20, 22).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/klagomål/A 52510-2022.docx", "A 52510-2022")'
$ws.Cells.Item(20, 23).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/klagomålsmail/A 52510-2022.docx", "A 52510-2022")'
$ws.Cells.Item(20, 24).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/tillsyn/A 52510-2022.docx", "A 52510-2022")'
$ws.Cells.Item(20, 25).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/tillsynsmail/A 52510-2022.docx", "A 52510-2022")'
$ws.Cells.Item(21, 19).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/artfynd/A 59665-2022.xlsx", "A 59665-2022")'
$ws.Cells.Item(21, 20).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/kartor/A 59665-2022.png", "A 59665-2022")'
$ws.Cells.Item(21, 21).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/knärot/A 59665-2022.png", "A 59665-2022")'
$ws.Cells.Item(21, 22).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/klagomål/A 59665-2022.docx", "A 59665-2022")'
$ws.Cells.Item(21, 23).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/klagomålsmail/A 59665-2022.docx", "A 59665-2022")'
$ws.Cells.Item(21, 24).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/tillsyn/A 59665-2022.docx", "A 59665-2022")'
$ws.Cells.Item(21, 25).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/tillsynsmail/A 59665-2022.docx", "A 59665-2022")'
$ws.Cells.Item(22, 19).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/artfynd/A 12314-2023.xlsx", "A 12314-2023")'
$ws.Cells.Item(22, 20).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/kartor/A 12314-2023.png", "A 12314-2023")'
$ws.Cells.Item(22, 22).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/klagomål/A 12314-2023.docx", "A 12314-2023")'
$ws.Cells.Item(22, 23).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/klagomålsmail/A 12314-2023.docx", "A 12314-2023")'
$ws.Cells.Item(22, 24).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/tillsyn/A 12314-2023.docx", "A 12314-2023")'
$ws.Cells.Item(22, 25).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/tillsynsmail/A 12314-2023.docx", "A 12314-2023")'
$ws.Cells.Item(23, 19).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/artfynd/A 12444-2023.xlsx", "A 12444-2023")'
$ws.Cells.Item(23, 20).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/kartor/A 12444-2023.png", "A 12444-2023")'
$ws.Cells.Item(23, 22).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/klagomål/A 12444-2023.docx", "A 12444-2023")'
$ws.Cells.Item(23, 23).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/klagomålsmail/A 12444-2023.docx", "A 12444-2023")'
$ws.Cells.Item(23, 24).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/tillsyn/A 12444-2023.docx", "A 12444-2023")'
$ws.Cells.Item(23, 25).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/tillsynsmail/A 12444-2023.docx", "A 12444-2023")'
$ws.Cells.Item(24, 19).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/artfynd/A 21926-2023.xlsx", "A 21926-2023")'
$ws.Cells.Item(24, 20).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/kartor/A 21926-2023.png", "A 21926-2023")'
$ws.Cells.Item(24, 21).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/knärot/A 21926-2023.png", "A 21926-2023")'
$ws.Cells.Item(24, 22).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/klagomål/A 21926-2023.docx", "A 21926-2023")'
$ws.Cells.Item(24, 23).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/klagomålsmail/A 21926-2023.docx", "A 21926-2023")'
$ws.Cells.Item(24, 24).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/tillsyn/A 21926-2023.docx", "A 21926-2023")'
$ws.Cells.Item(24, 25).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/tillsynsmail/A 21926-2023.docx", "A 21926-2023")'
$ws.Cells.Item(25, 19).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/artfynd/A 23652-2023.xlsx", "A 23652-2023")'
$ws.Cells.Item(25, 20).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/kartor/A 23652-2023.png", "A 23652-2023")'
$ws.Cells.Item(25, 22).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/klagomål/A 23652-2023.docx", "A 23652-2023")'
$ws.Cells.Item(25, 23).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/klagomålsmail/A 23652-2023.docx", "A 23652-2023")'
$ws.Cells.Item(25, 24).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/tillsyn/A 23652-2023.docx", "A 23652-2023")'
$ws.Cells.Item(25, 25).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/tillsynsmail/A 23652-2023.docx", "A 23652-2023")'
$ws.Cells.Item(26, 19).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/artfynd/A 31680-2023.xlsx", "A 31680-2023")'
$ws.Cells.Item(26, 20).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/kartor/A 31680-2023.png", "A 31680-2023")'
$ws.Cells.Item(26, 22).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/klagomål/A 31680-2023.docx", "A 31680-2023")'
$ws.Cells.Item(26, 23).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/klagomålsmail/A 31680-2023.docx", "A 31680-2023")'
$ws.Cells.Item(26, 24).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/tillsyn/A 31680-2023.docx", "A 31680-2023")'
$ws.Cells.Item(26, 25).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/tillsynsmail/A 31680-2023.docx", "A 31680-2023")'
$ws.Cells.Item(27, 19).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/artfynd/A 58399-2018.xlsx", "A 58399-2018")'
$ws.Cells.Item(27, 20).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/kartor/A 58399-2018.png", "A 58399-2018")'
$ws.Cells.Item(27, 22).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/klagomål/A 58399-2018.docx", "A 58399-2018")'
$ws.Cells.Item(27, 23).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/klagomålsmail/A 58399-2018.docx", "A 58399-2018")'
$ws.Cells.Item(27, 24).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/tillsyn/A 58399-2018.docx", "A 58399-2018")'
$ws.Cells.Item(27, 25).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/tillsynsmail/A 58399-2018.docx", "A 58399-2018")'
$ws.Cells.Item(28, 19).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/artfynd/A 66079-2018.xlsx", "A 66079-2018")'
$ws.Cells.Item(28, 20).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/kartor/A 66079-2018.png", "A 66079-2018")'
$ws.Cells.Item(28, 22).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/klagomål/A 66079-2018.docx", "A 66079-2018")'
$ws.Cells.Item(28, 23).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/klagomålsmail/A 66079-2018.docx", "A 66079-2018")'
$ws.Cells.Item(28, 24).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/tillsyn/A 66079-2018.docx", "A 66079-2018")'
$ws.Cells.Item(28, 25).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/tillsynsmail/A 66079-2018.docx", "A 66079-2018")'
$ws.Cells.Item(29, 19).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/artfynd/A 68272-2018.xlsx", "A 68272-2018")'
$ws.Cells.Item(29, 20).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/kartor/A 68272-2018.png", "A 68272-2018")'
$ws.Cells.Item(29, 22).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/klagomål/A 68272-2018.docx", "A 68272-2018")'
$ws.Cells.Item(29, 23).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/klagomålsmail/A 68272-2018.docx", "A 68272-2018")'
$ws.Cells.Item(29, 24).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/tillsyn/A 68272-2018.docx", "A 68272-2018")'
$ws.Cells.Item(29, 25).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/tillsynsmail/A 68272-2018.docx", "A 68272-2018")'
$ws.Cells.Item(30, 19).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/artfynd/A 16079-2019.xlsx", "A 16079-2019")'
$ws.Cells.Item(30, 20).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/kartor/A 16079-2019.png", "A 16079-2019")'
$ws.Cells.Item(30, 21).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/knärot/A 16079-2019.png", "A 16079-2019")'
$ws.Cells.Item(30, 22).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/klagomål/A 16079-2019.docx", "A 16079-2019")'
$ws.Cells.Item(30, 23).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/klagomålsmail/A 16079-2019.docx", "A 16079-2019")'
$ws.Cells.Item(30, 24).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/tillsyn/A 16079-2019.docx", "A 16079-2019")'
$ws.Cells.Item(30, 25).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/tillsynsmail/A 16079-2019.docx", "A 16079-2019")'
$ws.Cells.Item(31, 19).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/artfynd/A 23553-2020.xlsx", "A 23553-2020")'
$ws.Cells.Item(31, 20).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/kartor/A 23553-2020.png", "A 23553-2020")'
$ws.Cells.Item(31, 22).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/klagomål/A 23553-2020.docx", "A 23553-2020")'
$ws.Cells.Item(31, 23).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/klagomålsmail/A 23553-2020.docx", "A 23553-2020")'
$ws.Cells.Item(31, 24).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/tillsyn/A 23553-2020.docx", "A 23553-2020")'
$ws.Cells.Item(31, 25).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/tillsynsmail/A 23553-2020.docx", "A 23553-2020")'
$ws.Cells.Item(32, 19).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/artfynd/A 29908-2021.xlsx", "A 29908-2021")'
$ws.Cells.Item(32, 20).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/kartor/A 29908-2021.png", "A 29908-2021")'
$ws.Cells.Item(32, 22).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/klagomål/A 29908-2021.docx", "A 29908-2021")'
$ws.Cells.Item(32, 23).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/klagomålsmail/A 29908-2021.docx", "A 29908-2021")'
$ws.Cells.Item(32, 24).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/tillsyn/A 29908-2021.docx", "A 29908-2021")'
$ws.Cells.Item(32, 25).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/tillsynsmail/A 29908-2021.docx", "A 29908-2021")'
$ws.Cells.Item(33, 19).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/artfynd/A 38708-2021.xlsx", "A 38708-2021")'
$ws.Cells.Item(33, 20).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/kartor/A 38708-2021.png", "A 38708-2021")'
$ws.Cells.Item(33, 22).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/klagomål/A 38708-2021.docx", "A 38708-2021")'
$ws.Cells.Item(33, 23).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/klagomålsmail/A 38708-2021.docx", "A 38708-2021")'
$ws.Cells.Item(33, 24).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/tillsyn/A 38708-2021.docx", "A 38708-2021")'
$ws.Cells.Item(33, 25).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/tillsynsmail/A 38708-2021.docx", "A 38708-2021")'
$ws.Cells.Item(34, 19).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/artfynd/A 41874-2021.xlsx", "A 41874-2021")'
$ws.Cells.Item(34, 20).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/kartor/A 41874-2021.png", "A 41874-2021")'
$ws.Cells.Item(34, 22).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/klagomål/A 41874-2021.docx", "A 41874-2021")'
$ws.Cells.Item(34, 23).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/klagomålsmail/A 41874-2021.docx", "A 41874-2021")'
$ws.Cells.Item(34, 24).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/tillsyn/A 41874-2021.docx", "A 41874-2021")'
$ws.Cells.Item(34, 25).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/tillsynsmail/A 41874-2021.docx", "A 41874-2021")'
$ws.Cells.Item(35, 19).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/artfynd/A 49776-2021.xlsx", "A 49776-2021")'
$ws.Cells.Item(35, 20).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/kartor/A 49776-2021.png", "A 49776-2021")'
$ws.Cells.Item(35, 22).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/klagomål/A 49776-2021.docx", "A 49776-2021")'
$ws.Cells.Item(35, 23).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/klagomålsmail/A 49776-2021.docx", "A 49776-2021")'
$ws.Cells.Item(35, 24).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/tillsyn/A 49776-2021.docx", "A 49776-2021")'
$ws.Cells.Item(35, 25).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/tillsynsmail/A 49776-2021.docx", "A 49776-2021")'
$ws.Cells.Item(36, 19).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/artfynd/A 54045-2021.xlsx", "A 54045-2021")'
$ws.Cells.Item(36, 20).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/kartor/A 54045-2021.png", "A 54045-2021")'
$ws.Cells.Item(36, 22).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/klagomål/A 54045-2021.docx", "A 54045-2021")'
$ws.Cells.Item(36, 23).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/klagomålsmail/A 54045-2021.docx", "A 54045-2021")'
$ws.Cells.Item(36, 24).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/tillsyn/A 54045-2021.docx", "A 54045-2021")'
$ws.Cells.Item(36, 25).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/tillsynsmail/A 54045-2021.docx", "A 54045-2021")'
$ws.Cells.Item(37, 19).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/artfynd/A 54699-2021.xlsx", "A 54699-2021")'
$ws.Cells.Item(37, 20).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/kartor/A 54699-2021.png", "A 54699-2021")'
$ws.Cells.Item(37, 21).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/knärot/A 54699-2021.png", "A 54699-2021")'
$ws.Cells.Item(37, 22).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/klagomål/A 54699-2021.docx", "A 54699-2021")'
$ws.Cells.Item(37, 23).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/klagomålsmail/A 54699-2021.docx", "A 54699-2021")'
$ws.Cells.Item(37, 24).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/tillsyn/A 54699-2021.docx", "A 54699-2021")'
$ws.Cells.Item(37, 25).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/tillsynsmail/A 54699-2021.docx", "A 54699-2021")'
$ws.Cells.Item(38, 19).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/artfynd/A 6972-2022.xlsx", "A 6972-2022")'
$ws.Cells.Item(38, 20).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/kartor/A 6972-2022.png", "A 6972-2022")'
$ws.Cells.Item(38, 22).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/klagomål/A 6972-2022.docx", "A 6972-2022")'
$ws.Cells.Item(38, 23).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/klagomålsmail/A 6972-2022.docx", "A 6972-2022")'
$ws.Cells.Item(38, 24).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/tillsyn/A 6972-2022.docx", "A 6972-2022")'
$ws.Cells.Item(38, 25).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/tillsynsmail/A 6972-2022.docx", "A 6972-2022")'
$ws.Cells.Item(39, 19).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/artfynd/A 21302-2022.xlsx", "A 21302-2022")'
$ws.Cells.Item(39, 20).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/kartor/A 21302-2022.png", "A 21302-2022")'
$ws.Cells.Item(39, 22).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/klagomål/A 21302-2022.docx", "A 21302-2022")'
$ws.Cells.Item(39, 23).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/klagomålsmail/A 21302-2022.docx", "A 21302-2022")'
$ws.Cells.Item(39, 24).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/tillsyn/A 21302-2022.docx", "A 21302-2022")'
$ws.Cells.Item(39, 25).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/tillsynsmail/A 21302-2022.docx", "A 21302-2022")'
$ws.Cells.Item(40, 19).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/artfynd/A 35370-2022.xlsx", "A 35370-2022")'
$ws.Cells.Item(40, 20).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/kartor/A 35370-2022.png", "A 35370-2022")'
$ws.Cells.Item(40, 22).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/klagomål/A 35370-2022.docx", "A 35370-2022")'
$ws.Cells.Item(40, 23).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/klagomålsmail/A 35370-2022.docx", "A 35370-2022")'
$ws.Cells.Item(40, 24).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/tillsyn/A 35370-2022.docx", "A 35370-2022")'
$ws.Cells.Item(40, 25).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/tillsynsmail/A 35370-2022.docx", "A 35370-2022")'
$ws.Cells.Item(41, 19).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/artfynd/A 35389-2022.xlsx", "A 35389-2022")'
$ws.Cells.Item(41, 20).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/kartor/A 35389-2022.png", "A 35389-2022")'
$ws.Cells.Item(41, 22).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/klagomål/A 35389-2022.docx", "A 35389-2022")'
$ws.Cells.Item(41, 23).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/klagomålsmail/A 35389-2022.docx", "A 35389-2022")'
$ws.Cells.Item(41, 24).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/tillsyn/A 35389-2022.docx", "A 35389-2022")'
$ws.Cells.Item(41, 25).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/tillsynsmail/A 35389-2022.docx", "A 35389-2022")'
$ws.Cells.Item(42, 19).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/artfynd/A 57693-2022.xlsx", "A 57693-2022")'
$ws.Cells.Item(42, 20).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/kartor/A 57693-2022.png", "A 57693-2022")'
$ws.Cells.Item(42, 22).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/klagomål/A 57693-2022.docx", "A 57693-2022")'
$ws.Cells.Item(42, 23).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/klagomålsmail/A 57693-2022.docx", "A 57693-2022")'
$ws.Cells.Item(42, 24).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/tillsyn/A 57693-2022.docx", "A 57693-2022")'
$ws.Cells.Item(42, 25).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/tillsynsmail/A 57693-2022.docx", "A 57693-2022")'
$ws.Cells.Item(43, 19).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/artfynd/A 865-2023.xlsx", "A 865-2023")'
$ws.Cells.Item(43, 20).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/kartor/A 865-2023.png", "A 865-2023")'
$ws.Cells.Item(43, 21).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/knärot/A 865-2023.png", "A 865-2023")'
$ws.Cells.Item(43, 22).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/klagomål/A 865-2023.docx", "A 865-2023")'
$ws.Cells.Item(43, 23).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/klagomålsmail/A 865-2023.docx", "A 865-2023")'
$ws.Cells.Item(43, 24).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/tillsyn/A 865-2023.docx", "A 865-2023")'
$ws.Cells.Item(43, 25).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/tillsynsmail/A 865-2023.docx", "A 865-2023")'
$ws.Cells.Item(44, 19).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/artfynd/A 17571-2023.xlsx", "A 17571-2023")'
$ws.Cells.Item(44, 20).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/kartor/A 17571-2023.png", "A 17571-2023")'
$ws.Cells.Item(44, 22).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/klagomål/A 17571-2023.docx", "A 17571-2023")'
$ws.Cells.Item(44, 23).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/klagomålsmail/A 17571-2023.docx", "A 17571-2023")'
$ws.Cells.Item(44, 24).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/tillsyn/A 17571-2023.docx", "A 17571-2023")'
$ws.Cells.Item(44, 25).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/tillsynsmail/A 17571-2023.docx", "A 17571-2023")'
$ws.Cells.Item(45, 19).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/artfynd/A 17913-2023.xlsx", "A 17913-2023")'
$ws.Cells.Item(45, 20).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/kartor/A 17913-2023.png", "A 17913-2023")'
$ws.Cells.Item(45, 22).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/klagomål/A 17913-2023.docx", "A 17913-2023")'
$ws.Cells.Item(45, 23).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/klagomålsmail/A 17913-2023.docx", "A 17913-2023")'
$ws.Cells.Item(45, 24).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/tillsyn/A 17913-2023.docx", "A 17913-2023")'
$ws.Cells.Item(45, 25).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/tillsynsmail/A 17913-2023.docx", "A 17913-2023")'
$ws.Cells.Item(46, 19).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/artfynd/A 20925-2023.xlsx", "A 20925-2023")'
$ws.Cells.Item(46, 20).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/kartor/A 20925-2023.png", "A 20925-2023")'
$ws.Cells.Item(46, 21).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/knärot/A 20925-2023.png", "A 20925-2023")'
$ws.Cells.Item(46, 22).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/klagomål/A 20925-2023.docx", "A 20925-2023")'
$ws.Cells.Item(46, 23).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/klagomålsmail/A 20925-2023.docx", "A 20925-2023")'
$ws.Cells.Item(46, 24).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/tillsyn/A 20925-2023.docx", "A 20925-2023")'
$ws.Cells.Item(46, 25).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/tillsynsmail/A 20925-2023.docx", "A 20925-2023")'
$ws.Cells.Item(47, 19).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/artfynd/A 21131-2023.xlsx", "A 21131-2023")'
$ws.Cells.Item(47, 20).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/kartor/A 21131-2023.png", "A 21131-2023")'
$ws.Cells.Item(47, 22).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/klagomål/A 21131-2023.docx", "A 21131-2023")'
$ws.Cells.Item(47, 23).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/klagomålsmail/A 21131-2023.docx", "A 21131-2023")'
$ws.Cells.Item(47, 24).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/tillsyn/A 21131-2023.docx", "A 21131-2023")'
$ws.Cells.Item(47, 25).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/tillsynsmail/A 21131-2023.docx", "A 21131-2023")'
$ws.Cells.Item(48, 19).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/artfynd/A 22900-2023.xlsx", "A 22900-2023")'
$ws.Cells.Item(48, 20).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/kartor/A 22900-2023.png", "A 22900-2023")'
$ws.Cells.Item(48, 22).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/klagomål/A 22900-2023.docx", "A 22900-2023")'
$ws.Cells.Item(48, 23).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/klagomålsmail/A 22900-2023.docx", "A 22900-2023")'
$ws.Cells.Item(48, 24).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/tillsyn/A 22900-2023.docx", "A 22900-2023")'
$ws.Cells.Item(48, 25).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/tillsynsmail/A 22900-2023.docx", "A 22900-2023")'
$ws.Cells.Item(49, 19).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/artfynd/A 25763-2023.xlsx", "A 25763-2023")'
$ws.Cells.Item(49, 20).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/kartor/A 25763-2023.png", "A 25763-2023")'
$ws.Cells.Item(49, 22).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/klagomål/A 25763-2023.docx", "A 25763-2023")'
$ws.Cells.Item(49, 23).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/klagomålsmail/A 25763-2023.docx", "A 25763-2023")'
$ws.Cells.Item(49, 24).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/tillsyn/A 25763-2023.docx", "A 25763-2023")'
$ws.Cells.Item(49, 25).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/tillsynsmail/A 25763-2023.docx", "A 25763-2023")'
$ws.Cells.Item(182, 21).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/knärot/A 59674-2021.png", "A 59674-2021")'
$ws.Cells.Item(182, 22).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/klagomål/A 59674-2021.docx", "A 59674-2021")'
$ws.Cells.Item(182, 23).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/klagomålsmail/A 59674-2021.docx", "A 59674-2021")'
$ws.Cells.Item(182, 24).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/tillsyn/A 59674-2021.docx", "A 59674-2021")'
$ws.Cells.Item(182, 25).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/tillsynsmail/A 59674-2021.docx", "A 59674-2021")'
$ws.Cells.Item(256, 21).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/knärot/A 10751-2023.png", "A 10751-2023")'
$ws.Cells.Item(256, 22).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/klagomål/A 10751-2023.docx", "A 10751-2023")'
$ws.Cells.Item(256, 23).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/klagomålsmail/A 10751-2023.docx", "A 10751-2023")'
$ws.Cells.Item(256, 24).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/tillsyn/A 10751-2023.docx", "A 10751-2023")'
$ws.Cells.Item(256, 25).Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORLANGE/tillsynsmail/A 10751-2023.docx", "A 10751-2023")'
